$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.306.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.678.29'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.20%  '

$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.61%  '

$ws.Range("E6").Value = '  +0.77%  '

$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '30.16'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.264'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0622'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0902'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.920.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +17.77%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.690.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.86%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.617'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.24%  '

$ws.Range("E16").Value = '  +3.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.326.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.86'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '247.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0720'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.22%  '

$ws.Range("E21").Value = '  -0.22%  '

$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.62%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.98%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '

$ws.Range("E26").Value = '  +1.25%  '

$ws.Range("E28").Value = '  +2.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0501'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.14%  '

$ws.Range("E31").Value = '  +0.83%  '

$ws.Range("E32").Value = '  +3.75%  '

$ws.Range("E33").Value = '  +3.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.488.70'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.74'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0180'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.36%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.588'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.14%  '

$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '79.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.31'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.856'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.23%  '

$ws.Range("E43").Value = '  +2.16%  '

$ws.Range("E44").Value = '  +1.74%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.23%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '52.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.813.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.60%  '

$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '95.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0117'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.78%  '
